# ---------------------------------------------------------------------------
# Unidad_4/Presentación.pptx — "feat: update code and presentation"
#
# 1) Refresh the cached "datetimeFigureOut" field text (short date, e.g.
#    Mac-style "M/d/yy") everywhere it is shown: the slide master, every
#    slide layout, and the notes master.
# 2) Append a new slide 5 ("Usando Hooks") using the "Title and Content"
#    layout, with its title + bullet body text.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1) Update the cached date field text (9/18/2024 -> 9/18/24) ----------

$oldDate = "9/18/2024"
$newDate = "9/18/24"

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DateShapes $master.Shapes

# Every slide layout attached to the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShapes $layout.Shapes
}

# Notes master
$notesMaster = $p.NotesMaster
Update-DateShapes $notesMaster.Shapes

# --- 2) Append the new slide 5: "Usando Hooks" -----------------------------

$newIndex = $p.Slides.Count + 1
$slide = $p.Slides.Add($newIndex, 2)   # ppLayoutText/ppLayoutObject == "Title and Content"

$title = $slide.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Usando Hooks"

$body = $slide.Shapes.Item(2)
$body.TextFrame.TextRange.Text = "Crear un nuevo URL " + [char]0x201C + "users/form/:index?" + [char]0x201D + " para una vista de formulario.`rCrear un formulario para name, age, city y verified."
